$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.730.90'
$ws.Range('E2').Value = '  -2.21%  '
$ws.Range('D3').Value = '2.346.89'
$ws.Range('E3').Value = '  -3.23%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.32'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.69'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.636'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.89%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.616'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -8.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.12'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0921'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -3.06%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.42'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.38%  '
$ws.Range('E13').Value = '  -3.71%  '
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.98'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -8.15%  '
$ws.Range('D16').Value = '2.701.26'
$ws.Range('E16').Value = '  -3.24%  '
$ws.Range('D17').Value = '2.337.89'
$ws.Range('E17').Value = '  -3.82%  '
$ws.Range('D18').Value = '42.679.19'
$ws.Range('E18').Value = '  -2.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.67'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.64%  '
$ws.Range('E20').Value = '  -4.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '76.92'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('E22').Value = '  +2.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '260.07'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -5.24%  '
$ws.Range('E24').Value = '  -6.36%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.59'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.36'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -6.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.07'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '174.74'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.09'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -5.75%  '
$ws.Range('B32').Value = 'WEMIXToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.00'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -6.83%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.15'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0889'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -5.69%  '
$ws.Range('E35').Value = '  -2.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.117'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +6.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.61'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -6.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0360'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.79'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -7.01%  '
$ws.Range('E40').Value = '  -7.63%  '
$ws.Range('E41').Value = '  +1.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.47'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -10.39%  '
$ws.Range('E43').Value = '  -2.58%  '
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '115.47'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -10.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.60'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.85'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -7.63%  '
$ws.Range('E48').Value = '  -4.35%  '
$ws.Range('E49').Value = '  -6.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.40'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.95%  '
$ws.Range('E51').Value = '  -5.73%  '
